# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 59 (pushing existing rows 59-106 down to 60-107),
# carrying over the static/reference columns (Mercado, Region, Codreg, Tipo, Producto, Categoria,
# Calidad, Origen) from the old row 59, and filling in the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 59; everything below shifts down by one row.
$ws.Rows(59).Insert()

# Populate the newly inserted row 59 with the new weekly record.
$ws.Range("A59").Value = 11
$ws.Range("B59").Value = "Vega Monumental Concepción"
$ws.Range("C59").Value = "Bíobío"
$ws.Range("D59").Value = 44589
$ws.Range("D59").NumberFormat = $ws.Range("D60").NumberFormat
$ws.Range("E59").Value = 8
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100103
$ws.Range("H59").Value = "Frutos de hueso (carozo)"
$ws.Range("I59").Value = 100103001
$ws.Range("J59").Value = "Cereza"
$ws.Range("K59").Value = "Santina"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 160
$ws.Range("N59").Value = 7500
$ws.Range("O59").Value = 8000
$ws.Range("P59").Value = 7750
$ws.Range("Q59").Value = "$/bandeja 10 kilos"
$ws.Range("R59").Value = "Provincia de Curicó"
$ws.Range("S59").Value = 775
$ws.Range("T59").Value = 10
